# data : case 1
# Update the 32x2 numeric table (A1:B32) to the new computed values,
# and normalize both column widths to the same narrower width.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = -0.36509669890178031
$ws.Range("B1").Value = 0.36419359360247938
$ws.Range("A2").Value = -0.24274345885605975
$ws.Range("B2").Value = 0.24065892177893389
$ws.Range("A3").Value = -0.13771053154034263
$ws.Range("B3").Value = 0.13722874648107819
$ws.Range("A4").Value = -0.19922500794805487
$ws.Range("B4").Value = 0.19822643387503547
$ws.Range("A5").Value = -0.1922264342758826
$ws.Range("B5").Value = 0.19021163422745513
$ws.Range("A6").Value = -0.074340989448129058
$ws.Range("B6").Value = 0.074267089167018963
$ws.Range("A7").Value = -0.054267089659902012
$ws.Range("B7").Value = 0.054084472414890783
$ws.Range("A8").Value = -0.034084472911467323
$ws.Range("B8").Value = 0.033891332019213571
$ws.Range("A9").Value = -0.027891332444861305
$ws.Range("B9").Value = 0.027709788839300309
$ws.Range("A10").Value = -0.0075580034345605895
$ws.Range("B10").Value = 0.0075582366384168154
$ws.Range("A11").Value = -0.0030582370589833374
$ws.Range("B11").Value = 0.0030565516539304838
$ws.Range("A12").Value = 0.0029434479174401673
$ws.Range("B12").Value = -0.0029506768306593045
$ws.Range("A13").Value = 0.0089506764024154251
$ws.Range("B13").Value = -0.0089540777107197655
$ws.Range("A14").Value = 0.02095407725099907
$ws.Range("B14").Value = -0.020971344561718119
$ws.Range("A15").Value = 0.026971344134940622
$ws.Range("B15").Value = -0.027008085385835656
$ws.Range("A16").Value = 0.033008084960875372
$ws.Range("B16").Value = -0.033105374594412762
$ws.Range("A17").Value = 0.039105374173623808
$ws.Range("B17").Value = -0.039182502775344297
$ws.Range("A18").Value = -0.11145013514643054
$ws.Range("B18").Value = 0.11128840367103621
$ws.Range("A19").Value = -0.027097145993723704
$ws.Range("B19").Value = 0.027013748773502932
$ws.Range("A20").Value = -0.018013749173466209
$ws.Range("B20").Value = 0.018004303985692971
$ws.Range("A21").Value = -0.0090043043862406691
$ws.Range("B21").Value = 0.0089999995990508452
$ws.Range("A22").Value = -0.093948390644321478
$ws.Range("B22").Value = 0.093635084581233485
$ws.Range("A23").Value = -0.084635084990598131
$ws.Range("B23").Value = 0.084126937905486088
$ws.Range("A24").Value = -0.042126938497167643
$ws.Range("B24").Value = 0.041999999404983335
$ws.Range("A25").Value = -0.094943253012427675
$ws.Range("B25").Value = 0.09469756434834764
$ws.Range("A26").Value = -0.088697564767233672
$ws.Range("B26").Value = 0.088380766957772749
$ws.Range("A27").Value = -0.082380767379118147
$ws.Range("B27").Value = 0.081296591040773691
$ws.Range("A28").Value = -0.075296591471190943
$ws.Range("B28").Value = 0.074545368421570224
$ws.Range("A29").Value = -0.062545368890431163
$ws.Range("B29").Value = 0.062172916003341072
$ws.Range("A30").Value = -0.042172916517877912
$ws.Range("B30").Value = 0.042020230301042183
$ws.Range("A31").Value = -0.027020230793461053
$ws.Range("B31").Value = 0.027000809856669861
$ws.Range("A32").Value = -0.0060008103816944214
$ws.Range("B32").Value = 0.0059999995545672036

# Column A was 15.7109375 and column B was 16.42578125; both become 15.42578125.
# ColumnWidth is stored by Excel on a whole-pixel grid, so we feed it the raw
# character-width value that rounds to the closest achievable stored width.
$ws.Columns.Item(1).ColumnWidth = 14.6
$ws.Columns.Item(2).ColumnWidth = 14.6
